{"js": "// Replace the five cover-letter body sentences with their revised wording.\n// Each original sentence lives in its own `w:r` run inside one paragraph\n// (the runs are separated by `w:br` line-break runs), so a body-wide\n// exact-text search finds a single, unambiguous hit per sentence that we\n// then swap in place with Range.insertText(..., \"Replace\").\nconst replacements = [\n  [\n    \"I am writing to express my interest in the position you have available. Although the job title, company, and description were not specified, I am eager to contribute my skills and experience to your team and help achieve your organizational goals.\",\n    \"I am excited to apply for the position at your company, where I hope to contribute my skills and enthusiasm. Although the job title and company details were unspecified, I am eager to bring my dedication and adaptability to any role presented.\"\n  ],\n  [\n    \"My background includes working with various technologies and adapting quickly to new environments. I am confident that my ability to learn rapidly and apply knowledge effectively will be an asset to your company.\",\n    \"My background includes experience with various technologies and a commitment to continuous learning. I am confident that my ability to quickly grasp new concepts and technologies will allow me to effectively support your team's goals.\"\n  ],\n  [\n    \"I am particularly interested in roles that offer the opportunity to grow professionally and tackle challenging projects. I am enthusiastic about collaborating with colleagues to develop innovative solutions and improve existing processes.\",\n    \"I thrive in collaborative environments and enjoy working with diverse teams to solve complex problems. I am particularly drawn to opportunities that challenge me and allow for professional growth.\"\n  ],\n  [\n    \"Throughout my career, I have demonstrated strong communication skills and a commitment to delivering high-quality work. I am confident that these qualities will enable me to make a positive impact in your organization.\",\n    \"I welcome the chance to discuss how my skills and experiences could be a match for your needs. Please feel free to contact me at your convenience to arrange a conversation.\"\n  ],\n  [\n    \"Thank you for considering my application. I look forward to the opportunity to discuss how my skills and experience align with your needs and contribute to your company's success.\",\n    \"Thank you for considering my application. I look forward to the possibility of contributing to your organization and am enthusiastic about the potential to grow together.\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find expected text: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the five cover-letter body sentences with their revised wording.\n# Each original sentence lives in its own run inside one paragraph (the runs\n# are separated by manual line-break runs), so Find.Execute against the\n# whole-document range locates a single, unambiguous match per sentence.\n# Assigning the matched Range's .Text directly (rather than passing the new\n# wording through Find.Replacement.Text) avoids Word's \"smart quotes\"\n# AutoFormat-as-you-type substitution, which would otherwise turn the plain\n# apostrophe in \"team's\" into a curly one.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"I am writing to express my interest in the position you have available. Although the job title, company, and description were not specified, I am eager to contribute my skills and experience to your team and help achieve your organizational goals.\"\n        New = \"I am excited to apply for the position at your company, where I hope to contribute my skills and enthusiasm. Although the job title and company details were unspecified, I am eager to bring my dedication and adaptability to any role presented.\"\n    },\n    @{\n        Old = \"My background includes working with various technologies and adapting quickly to new environments. I am confident that my ability to learn rapidly and apply knowledge effectively will be an asset to your company.\"\n        New = \"My background includes experience with various technologies and a commitment to continuous learning. I am confident that my ability to quickly grasp new concepts and technologies will allow me to effectively support your team's goals.\"\n    },\n    @{\n        Old = \"I am particularly interested in roles that offer the opportunity to grow professionally and tackle challenging projects. I am enthusiastic about collaborating with colleagues to develop innovative solutions and improve existing processes.\"\n        New = \"I thrive in collaborative environments and enjoy working with diverse teams to solve complex problems. I am particularly drawn to opportunities that challenge me and allow for professional growth.\"\n    },\n    @{\n        Old = \"Throughout my career, I have demonstrated strong communication skills and a commitment to delivering high-quality work. I am confident that these qualities will enable me to make a positive impact in your organization.\"\n        New = \"I welcome the chance to discuss how my skills and experiences could be a match for your needs. Please feel free to contact me at your convenience to arrange a conversation.\"\n    },\n    @{\n        Old = \"Thank you for considering my application. I look forward to the opportunity to discuss how my skills and experience align with your needs and contribute to your company's success.\"\n        New = \"Thank you for considering my application. I look forward to the possibility of contributing to your organization and am enthusiastic about the potential to grow together.\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $found = $rng.Find.Execute($r.Old)\n    if (-not $found) {\n        throw \"Could not find expected text: $($r.Old)\"\n    }\n    $rng.Text = $r.New\n}\n"}
